$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old extra model rows (rows 6-14); only 4 data rows remain (rows 2-5)
$ws.Rows("6:14").Delete()

# Overwrite the remaining 4 data rows with the updated AIC-based model comparison
$ws.Range("A2").Value = "Species"
$ws.Range("B2").Value = -553.6
$ws.Range("C2").Value = 0.8618
$ws.Range("D2").Value = 281.8

$ws.Range("A3").Value = "Species + Sex"
$ws.Range("B3").Value = -549.70000000000005
$ws.Range("C3").Value = 0.86250000000000004
$ws.Range("D3").Value = 280.8

$ws.Range("A4").Value = "Species * Sex"
$ws.Range("B4").Value = -545.1
$ws.Range("C4").Value = 0.86350000000000005
$ws.Range("D4").Value = 281.5

$ws.Range("A5").Value = "Sex"
$ws.Range("B5").Value = -160.9
$ws.Range("C5").Value = 0.0141
$ws.Range("D5").Value = 83.4

# Re-apply the sort so the sortState/autofilter-ish metadata shrinks to the new range
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B5"))
$ws.Sort.SetRange($ws.Range("A1:D5"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Match the reviewer's zoomed-in view and full-table selection
$excel.ActiveWindow.Zoom = 160
$ws.Range("A1:D5").Select()
